$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix typos: add a missing space before the line-break marker (or, for A23,
# replace the line-break marker with a plain space) in the fortune texts.
# (Order matches the original authoring order so shared-string indices line up.)
$ws.Range("A23").Value = "당신이 믿는 신념은 틀리지 않을 거에요. 이루미도 그것을 믿어요"
$ws.Range("A9").Value  = "감정을 숨기지 마세요. \r\n상대는 당신의 진심을 기다리고 있어요."
$ws.Range("A10").Value = "지금 하고 싶은 일이 있군요. \r\n망설이지 마세요 잘 될 거에요"
$ws.Range("A11").Value = "오늘 행운의 색은 푸른색입니다. \r\n시원시원한 모습을 보여주세요."
$ws.Range("A12").Value = "작은 말이라도 칭찬해보세요. \r\n고래가 춤을 출지도 모르니까요."
$ws.Range("A14").Value = "너무 과한 계획을 세우지는 않았나요. \r\n재정비를 해보세요."
$ws.Range("A16").Value = "영화나 연극을 보세요. \r\n생각지 못한 영감을 얻을 것 같아요."
$ws.Range("A17").Value = "당분간 지켜보는 것이 좋을 것 같아요. \r\n타이밍이 중요해요"
$ws.Range("A18").Value = "노력은 배신하지 않아요. \r\n성실한 당신을 이루미는 믿어요."
$ws.Range("A20").Value = "당신은 웃는 얼굴이 참 예뻐요. \r\n당신의 미소가 힘이 돼요."
$ws.Range("A21").Value = "할 수 있다고 믿어 봐요. \r\n생각보다 별 것 아닐 수도 있어요."
$ws.Range("A24").Value = "일 년 전 당신을 떠올려 보세요. \r\n분명 계속 성장했어요."
$ws.Range("A31").Value = "우리는 아직 우리가 보석인지 몰라요. \r\n당신은 특별해요"

# Update the sheet's stored selection to reflect where the author was last
# working in the file (A13 -> A31).
$ws.Range("A31").Select()
